# Add a new "network" entry to the command library on sheet 1 (row 4)
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("A4").Value = "network"
$ws.Range("B4").Value = "make a network share drive"
$ws.Range("C4").Value = "1. Folder right click > properties > sharing tab > share > add everyone to write/read`n2. >net use a: \\pc-011-032\drivea`n3. >net use     **this can check whether the new netowrk shared folder is activated"

# Wrap the long instructions text and size the row to fit it
$ws.Range("C4").WrapText = $true
$ws.Rows.Item(4).RowHeight = 47.25

# Move the active selection the way the author left it
$ws.Range("C5").Select() | Out-Null
